# Automatische test-sync: 2025-07-29 21:51:50
# Adds a new "Testmail #11" (retour niet verwerkt) log row to the historical
# responses sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 24

$ws.Cells.Item($newRow, 1).Value  = "Testmail #11: Mijn retour is nog steeds niet verwerkt."
$ws.Cells.Item($newRow, 2).Value  = "Beste klant,`nBedankt voor je bericht. Het spijt me te horen dat je retour nog niet is verwerkt. Om dit verder te kunnen onderzoeken, heb ik wat meer informatie van je nodig. Zou je alsjeblieft je ordernummer en/of het trackingnummer van de retourzending kunnen doorgeven? Hiermee kunnen we precies nakijken wat er gaande is en je zo snel mogelijk helpen.`nMet vriendelijke groet,`n[Naam]`nKlantenservice Team"
$ws.Cells.Item($newRow, 3).Value  = "Mijn retour is nog steeds niet verwerkt."
$ws.Cells.Item($newRow, 4).Value  = "mailmind.test@zohomail.eu"
$ws.Cells.Item($newRow, 5).Value  = "Retour / Terugbetaling"
$ws.Cells.Item($newRow, 6).Value  = "2025-07-29 21:51:07"
$ws.Cells.Item($newRow, 7).Value  = "Ja"
$ws.Cells.Item($newRow, 8).Value  = "Nee"
$ws.Cells.Item($newRow, 9).Value  = "Ja"
$ws.Cells.Item($newRow, 10).Value = "Nee"
